# permohonan + inv barang
# Insert a new "Total Harga (Rp.)" column into the request-items table,
# between "Jumlah Barang" and "Keterangan", and fill in its header /
# placeholder row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift columns F (Keterangan) and G (Status) one to the right, opening up
# a blank column F for the new "Total Harga (Rp.)" field.
$ws.Columns.Item(6).Insert() | Out-Null

# New column's width (xlsx stored width 33 == ColumnWidth 32.1666...).
$ws.Columns.Item(6).ColumnWidth = 32.16666666666667

# Header + placeholder text for the new column.
$ws.Range("F6").Value = "Total Harga (Rp.)"
$ws.Range("F7").Value = "[a.totalharga]"

# Grow the table (ListObject) to cover the new column, then restore the
# names of the two columns that were pushed right (the resize operation
# re-derives names from header text and can mis-assign them).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B6:H7")) | Out-Null
$ws.Range("G6").Value = "Keterangan"
$ws.Range("H6").Value = "Status"

# Match the recorded selection state.
$ws.Range("F12").Select() | Out-Null
